$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.77539444320486
$ws.Range("C2").Value = 5.618615777951677
$ws.Range("E2").Value = 12.75368349509711
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 3.644333388148898
$ws.Range("K2").Value = 8.016937088703637
$ws.Range("M2").Value = 13.44571191124559
$ws.Range("O2").Value = 22.85787117799837
$ws.Range("B3").Value = 8.493080847653289
$ws.Range("C3").Value = 5.523668426412742
$ws.Range("E3").Value = 12.54331243916406
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 3.646118907752443
$ws.Range("K3").Value = 7.825180600948357
$ws.Range("M3").Value = 13.27559894548341
$ws.Range("O3").Value = 22.96939555485234
$ws.Range("B4").Value = 8.316184986892294
$ws.Range("C4").Value = 5.464116026518473
$ws.Range("E4").Value = 12.416936497438
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 3.647272281262355
$ws.Range("K4").Value = 7.706038463641009
$ws.Range("M4").Value = 13.17319416131309
$ws.Range("O4").Value = 23.04316136907838
$ws.Range("B5").Value = 8.243317892729145
$ws.Range("C5").Value = 5.439551981554254
$ws.Range("E5").Value = 12.3662138255407
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 3.647756683548254
$ws.Range("K5").Value = 7.65720281892831
$ws.Range("M5").Value = 13.13202642334059
$ws.Range("O5").Value = 23.07454882229215
$ws.Range("B6").Value = 8.231174702463518
$ws.Range("C6").Value = 5.435455847324894
$ws.Range("E6").Value = 12.35784040855884
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 3.647837988837925
$ws.Range("K6").Value = 7.649078590587795
$ws.Range("M6").Value = 13.12522596163587
$ws.Range("O6").Value = 23.07984076933399
$ws.Range("B7").Value = 8.315205280318786
$ws.Range("C7").Value = 5.463785919332838
$ws.Range("E7").Value = 12.41624919242777
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 3.647278755743048
$ws.Range("K7").Value = 7.705380907241127
$ws.Range("M7").Value = 13.17263661545858
$ws.Range("O7").Value = 23.04357930079218
$ws.Range("B8").Value = 8.678855473152504
$ws.Range("C8").Value = 5.586149723015797
$ws.Range("E8").Value = 12.68061288566543
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 3.644937220648405
$ws.Range("K8").Value = 7.951150746717353
$ws.Range("M8").Value = 13.38666201014062
$ws.Range("O8").Value = 22.89522529335342
$ws.Range("B9").Value = 9.359260990325335
$ws.Range("C9").Value = 5.815335401463749
$ws.Range("E9").Value = 13.21787462590273
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.640796109785169
$ws.Range("K9").Value = 8.419130156083522
$ws.Range("M9").Value = 13.82041054280421
$ws.Range("O9").Value = 22.64639739797254
$ws.Range("B10").Value = 9.833815482209079
$ws.Range("C10").Value = 5.976116478264029
$ws.Range("E10").Value = 13.6196430099036
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.638025417467985
$ws.Range("K10").Value = 8.750884649862096
$ws.Range("M10").Value = 14.14469507608299
$ws.Range("O10").Value = 22.48941723614329
$ws.Range("B11").Value = 10.04323777419087
$ws.Range("C11").Value = 6.047406830219406
$ws.Range("E11").Value = 13.8030210983532
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.636823347400271
$ws.Range("K11").Value = 8.898514504795205
$ws.Range("M11").Value = 14.29281312427323
$ws.Range("O11").Value = 22.42365034881064
$ws.Range("B12").Value = 10.12154344621468
$ws.Range("C12").Value = 6.074120864766599
$ws.Range("E12").Value = 13.87247483952577
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.636376495500699
$ws.Range("K12").Value = 8.953894828762573
$ws.Range("M12").Value = 14.3489362985568
$ws.Range("O12").Value = 22.39956107619441
$ws.Range("B13").Value = 10.10472433334886
$ws.Range("C13").Value = 6.068380310911335
$ws.Range("E13").Value = 13.85751735507023
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.636472362531048
$ws.Range("K13").Value = 8.941991757976506
$ws.Range("M13").Value = 14.33684847096796
$ws.Range("O13").Value = 22.40471281723808
$ws.Range("B14").Value = 10.04970040236767
$ws.Range("C14").Value = 6.049610351903795
$ws.Range("E14").Value = 13.808735273456
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.636786417600479
$ws.Range("K14").Value = 8.903081456020077
$ws.Range("M14").Value = 14.29743002312015
$ws.Range("O14").Value = 22.42165215088973
$ws.Range("B15").Value = 10.01586474893462
$ws.Range("C15").Value = 6.038076014280594
$ws.Range("E15").Value = 13.778854262642
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.63697987104335
$ws.Range("K15").Value = 8.879178094746383
$ws.Range("M15").Value = 14.27328799617215
$ws.Range("O15").Value = 22.43213425500254
$ws.Range("B16").Value = 9.819993206113116
$ws.Range("C16").Value = 5.971418868230895
$ws.Range("E16").Value = 13.60766461577673
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.63810514575313
$ws.Range("K16").Value = 8.741166076233579
$ws.Range("M16").Value = 14.13502279672171
$ws.Range("O16").Value = 22.49382913847063
$ws.Range("B17").Value = 9.698126905061567
$ws.Range("C17").Value = 5.930041362402711
$ws.Range("E17").Value = 13.50274745880875
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 3.638810375455469
$ws.Range("K17").Value = 8.655620548965238
$ws.Range("M17").Value = 14.05031770241446
$ws.Range("O17").Value = 22.53312510429637
$ws.Range("B18").Value = 9.627428714019844
$ws.Range("C18").Value = 5.906068909599393
$ws.Range("E18").Value = 13.44246383030706
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 3.639221497885159
$ws.Range("K18").Value = 8.606110444882182
$ws.Range("M18").Value = 14.00165705216032
$ws.Range("O18").Value = 22.55625814325384
$ws.Range("B19").Value = 9.603390054765914
$ws.Range("C19").Value = 5.897923005909559
$ws.Range("E19").Value = 13.42206581024405
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 3.639361641584568
$ws.Range("K19").Value = 8.589296147672391
$ws.Range("M19").Value = 13.98519321205573
$ws.Range("O19").Value = 22.56418169038167
$ws.Range("B20").Value = 9.711162814722412
$ws.Range("C20").Value = 5.934464114129446
$ws.Range("E20").Value = 13.51391021323481
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 3.638734734354812
$ws.Range("K20").Value = 8.664759130393408
$ws.Range("M20").Value = 14.05932892914522
$ws.Range("O20").Value = 22.52888699519691
$ws.Range("B21").Value = 10.06588986761839
$ws.Range("C21").Value = 6.055131322804334
$ws.Range("E21").Value = 13.82306400467791
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.636693945840674
$ws.Range("K21").Value = 8.914524961405313
$ws.Range("M21").Value = 14.3090076616407
$ws.Range("O21").Value = 22.41665450428078
$ws.Range("B22").Value = 10.2918804125497
$ws.Range("C22").Value = 6.132343372971304
$ws.Range("E22").Value = 14.02515327562556
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.635408799827414
$ws.Range("K22").Value = 9.094088621575573
$ws.Range("M22").Value = 14.47236246582668
$ws.Range("O22").Value = 22.34805750807498
$ws.Range("B23").Value = 10.17182079661348
$ws.Range("C23").Value = 6.091289981549657
$ws.Range("E23").Value = 13.91731508080605
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.636090270869921
$ws.Range("K23").Value = 8.989502608644489
$ws.Range("M23").Value = 14.38517781281172
$ws.Range("O23").Value = 22.38423287902963
$ws.Range("B24").Value = 9.705271254430208
$ws.Range("C24").Value = 5.932465159902944
$ws.Range("E24").Value = 13.50886341828592
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.638768914010754
$ws.Range("K24").Value = 8.660628598054137
$ws.Range("M24").Value = 14.05525483293919
$ws.Range("O24").Value = 22.53080135854156
$ws.Range("B25").Value = 9.179287653242534
$ws.Range("C25").Value = 5.754597782649912
$ws.Range("E25").Value = 13.07097493478446
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.641868451337158
$ws.Range("K25").Value = 8.294407015216938
$ws.Range("M25").Value = 14.27328799617215
$ws.Range("O25").Value = 22.43213425500254
